# Actualizar menu desde Excel
# Add a new "Fiscalia" district row to the Deliveries sheet's table (Tabla3),
# then make the Deliveries sheet the active tab with the new cell selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Deliveries")
$tbl = $ws.ListObjects.Item("Tabla3")

# Grow the table by one row (this also extends the table ref/autoFilter and
# the sheet dimension automatically).
$newRow = $tbl.ListRows.Add()

$ws.Range("A6").Value = "Fiscalia"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = "1 pm llega el menu"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0

# Make Deliveries the active sheet/tab, with C6 selected.
$ws.Activate()
$ws.Range("C6").Select()
